$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 346, shifting existing rows 346-406 down to 348-408
$ws.Rows.Item(346).Insert()
$ws.Rows.Item(346).Insert()

# Fill row 346 (new record)
$ws.Range("A346").Value = 5
$ws.Range("B346").Value = "Macroferia Regional de Talca"
$ws.Range("C346").Value = "Maule"
$ws.Range("D346").Value = 44504
$ws.Range("D346").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E346").Value = 7
$ws.Range("F346").Value = "Fruta"
$ws.Range("G346").Value = 100108
$ws.Range("H346").Value = "Tropicales y subtropicales"
$ws.Range("I346").Value = 100108006
$ws.Range("J346").Value = "Plátano"
$ws.Range("K346").Value = "Sin especificar"
$ws.Range("L346").Value = "Pintón"
$ws.Range("M346").Value = 480
$ws.Range("N346").Value = 15000
$ws.Range("O346").Value = 15000
$ws.Range("P346").Value = 15000
$ws.Range("Q346").Value = "`$/caja 20 kilos"
$ws.Range("R346").Value = "Ecuador"
$ws.Range("S346").Value = 750
$ws.Range("T346").Value = 20

# Fill row 347 (new record)
$ws.Range("A347").Value = 5
$ws.Range("B347").Value = "Macroferia Regional de Talca"
$ws.Range("C347").Value = "Maule"
$ws.Range("D347").Value = 44504
$ws.Range("D347").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E347").Value = 7
$ws.Range("F347").Value = "Fruta"
$ws.Range("G347").Value = 100108
$ws.Range("H347").Value = "Tropicales y subtropicales"
$ws.Range("I347").Value = 100108006
$ws.Range("J347").Value = "Plátano"
$ws.Range("K347").Value = "Sin especificar"
$ws.Range("L347").Value = "Primera Pintón"
$ws.Range("M347").Value = 360
$ws.Range("N347").Value = 16000
$ws.Range("O347").Value = 16000
$ws.Range("P347").Value = 16000
$ws.Range("Q347").Value = "`$/caja 20 kilos"
$ws.Range("R347").Value = "Ecuador"
$ws.Range("S347").Value = 800
$ws.Range("T347").Value = 20

Write-Host "Done"
